{"js": "// Split the single run of text in a few specific paragraphs into one run\n// per \"word\" and one run per run-of-whitespace, e.g.\n//   \"New guide!\" -> [\"New\"] [\" \"] [\"guide!\"]\n// This mirrors the target diff, which breaks a single <w:r> into many\n// <w:r> siblings (one per token), without altering the paragraph's\n// formatting (pPr) or the visible text itself.\n\nconst targets = new Set([\n  \"New guide!\",\n  \"Tom Coleman\",\n  \"Guide on rationalizing the denominator available now!\",\n]);\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Tokenize into alternating runs of non-whitespace / whitespace, e.g.\n// \"Guide on rationalizing\" -> [\"Guide\", \" \", \"on\", \" \", \"rationalizing\"].\nfunction tokenize(text) {\n  return text.match(/\\S+|\\s+/g) || [text];\n}\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\");\n}\n\n// Build the flat-OPC \"WordOpenXML\" wrapper insertOoxml() requires, a single\n// <w:p> with the original <w:pPr> (if any) and one <w:r><w:t>...</w:t></w:r>\n// per token.\nfunction buildParagraphOoxml(pPrXml, tokens) {\n  const runs = tokens\n    .map((t) => `<w:r><w:t xml:space=\"preserve\">${escapeXml(t)}</w:t></w:r>`)\n    .join(\"\");\n  return `<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body><w:p>${pPrXml || \"\"}${runs}</w:p></w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n}\n\nconst matches = paragraphs.items.filter((p) => targets.has(p.text));\n\nfor (const p of matches) {\n  // Grab this paragraph's own OOXML so we can carry its <w:pPr> (style,\n  // numbering, etc.) over unchanged into the rebuilt paragraph.\n  const ooxmlResult = p.getOoxml();\n  await context.sync();\n\n  const xml = ooxmlResult.value;\n  const pPrMatch = xml.match(/<w:pPr\\b[^>]*>[\\s\\S]*?<\\/w:pPr>|<w:pPr\\b[^>]*\\/>/);\n  const pPrXml = pPrMatch ? pPrMatch[0] : \"\";\n\n  const tokens = tokenize(p.text);\n  const flatOpc = buildParagraphOoxml(pPrXml, tokens);\n\n  p.getRange().insertOoxml(flatOpc, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Split the single run of text in a few specific paragraphs into one run\n# per \"word\" and one run per run-of-whitespace, e.g.\n#   \"New guide!\" -> [\"New\"] [\" \"] [\"guide!\"]\n# This mirrors the target diff, which breaks a single <w:r> into many\n# <w:r> siblings (one per token), without altering the paragraph's\n# formatting (pPr) or the visible text itself.\n\nfunction Escape-Xml($s) {\n    $s = $s -replace '&', '&amp;'\n    $s = $s -replace '<', '&lt;'\n    $s = $s -replace '>', '&gt;'\n    $s = $s -replace '\"', '&quot;'\n    return $s\n}\n\nfunction Tokenize($text) {\n    # Alternating runs of non-whitespace / whitespace, e.g.\n    # \"Guide on rationalizing\" -> [\"Guide\", \" \", \"on\", \" \", \"rationalizing\"].\n    $toks = @()\n    foreach ($m in [regex]::Matches($text, '\\S+|\\s+')) {\n        $toks += $m.Value\n    }\n    return $toks\n}\n\nfunction Build-ParagraphOoxml($pPrXml, $tokens) {\n    $runs = \"\"\n    foreach ($t in $tokens) {\n        $runs += \"<w:r><w:t xml:space=`\"preserve`\">$(Escape-Xml $t)</w:t></w:r>\"\n    }\n    return @\"\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body><w:p>$pPrXml$runs</w:p></w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n\"@\n}\n\n$targets = @(\n    \"New guide!\",\n    \"Tom Coleman\",\n    \"Guide on rationalizing the denominator available now!\"\n)\n\n$d = $word.ActiveDocument\n\n# Snapshot paragraphs first since rewriting one paragraph's XML can\n# invalidate/reindex the live COM collection while we are iterating it.\n$toEdit = @()\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    # Paragraph.Range.Text includes the trailing paragraph-mark char; trim it\n    # before comparing against the target plain-text strings.\n    $plain = $text.TrimEnd(\"`r\")\n    if ($targets -contains $plain) {\n        $toEdit += $p.Range\n    }\n}\n\nforeach ($range in $toEdit) {\n    $text = $range.Text.TrimEnd(\"`r\")\n    $xml = $range.WordOpenXML\n    $m = [regex]::Match($xml, '<w:pPr\\b[^>]*>[\\s\\S]*?</w:pPr>|<w:pPr\\b[^>]*/>')\n    $pPrXml = \"\"\n    if ($m.Success) {\n        $pPrXml = $m.Value\n    }\n\n    $tokens = Tokenize($text)\n    $flatOpc = Build-ParagraphOoxml $pPrXml $tokens\n\n    $range.InsertXML($flatOpc)\n}\n"}
